# "People can move now. Will register live on other's screens"
#
# The sheet holds a checkerboard-style grid of 0/1 values (testmap2 /
# testing-dnd). A piece has moved: the two cells that used to read
# 0/1 at C2:D2 and 1/0 at C3:D3 are now both marked "2" (the piece's
# new occupied squares), and the live selection follows the piece to
# its new square (C3) so collaborators watching the sheet see where it
# landed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2

# Move / broadcast the active selection to the piece's new square.
$ws.Range("C3").Select()
